$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update MACRO_SCORE column (N) for rows 2-6 (same new value across all rows)
$ws.Range("N2").Value = 85.92117485762657
$ws.Range("N3").Value = 85.92117485762657
$ws.Range("N4").Value = 85.92117485762657
$ws.Range("N5").Value = 85.92117485762657
$ws.Range("N6").Value = 85.92117485762657

# Update Bitcoin (BTC-USD) row values (row 4)
$ws.Range("D4").Value = 90977.63
$ws.Range("F4").Value = 4.81
$ws.Range("H4").Value = 56
$ws.Range("I4").Value = 43
$ws.Range("J4").Value = 43
$ws.Range("K4").Value = 52
